# Daily attendance processing - 2025-12-07 17:49:12
# Reorders the "Recorded By" (column G) value so that any leading
# "System"/"system" token(s) are moved to the end of the comma-separated list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -and $val.ToString().StartsWith("System, ")) {
        $parts = $val.ToString().Split(",")
        for ($i = 0; $i -lt $parts.Length; $i++) {
            $parts[$i] = $parts[$i].Trim()
        }

        $i = 0
        while ($i -lt $parts.Length -and ($parts[$i] -eq "System" -or $parts[$i] -eq "system")) {
            $i++
        }

        $leading = $parts[0..($i - 1)]
        $rest = @()
        if ($i -lt $parts.Length) {
            $rest = $parts[$i..($parts.Length - 1)]
        }

        $newParts = $rest + $leading
        $newVal = [string]::Join(", ", $newParts)

        $cell.Value = $newVal
    }
}
